$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (DATE) holds the same date serial value (45498 = 2024-07-25) for
# every data row from 4 to 151. The commit bumps that date forward one week
# to 45505 (2024-08-01) for all of those rows.
for ($row = 4; $row -le 151; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq 45498) {
        $cell.Value2 = 45505
    }
}
